$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 93

# The "Date" column in this sheet stores dates as plain text (e.g. "11/17/2025"),
# not as real date serials. Assigning a date-shaped string directly would make
# Excel auto-convert it to a date value, so we briefly force the cell to a
# text format while writing the value, then clear the formatting again so the
# new cell keeps the sheet's default (unstyled) look, matching the rows above it.
$dateCell = $ws.Range("A$row")
$dateCell.NumberFormat = "@"
$dateCell.Value = "11/18/2025"
$dateCell.ClearFormats()

$ws.Range("B$row").Value = 8123.83
